$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计").
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $afterSheet)
$q1.Name = "2022-Q1"

# Header row formatting (bold, thin border, centered / top aligned) to match
# the other fund-holding sheets.
$headerRange = $q1.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Text-valued data columns (force text so things like "008763" / "92.10"
# keep their leading/trailing zeros instead of becoming numbers).
$q1.Range("B2:G3").NumberFormat = "@"

$idxRange = $q1.Range("A2:A3")
$idxRange.HorizontalAlignment = -4108
$idxRange.VerticalAlignment = -4160
$idxRange.Borders.LineStyle = 1
$idxRange.Font.Bold = $true

$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "008763"
$q1.Range("C2").Value = "天弘越南市场股票（QDII）A"
$q1.Range("D2").Value = "37.53"
$q1.Range("E2").Value = "92.10"
$q1.Range("F2").Value = "4.57"
$q1.Range("G2").Value = "1.7151"
$q1.Range("H2").Value = 8

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "008764"
$q1.Range("C3").Value = "天弘越南市场股票（QDII）C"
$q1.Range("D3").Value = "14.26"
$q1.Range("E3").Value = "92.10"
$q1.Range("F3").Value = "4.57"
$q1.Range("G3").Value = "0.6517"
$q1.Range("H3").Value = 8

# ---------------------------------------------------------------------------
# 2. Prepend a 2022-Q1 summary row to the "总计" sheet, pushing the existing
#    rows down (2021-Q4 -> row 3, 2021-Q2 -> row 4).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$oldB2 = $total.Range("B2").Text
$oldC2 = $total.Range("C2").Value2
$oldD2 = $total.Range("D2").Value2
$oldB3 = $total.Range("B3").Text
$oldC3 = $total.Range("C3").Value2
$oldD3 = $total.Range("D3").Value2

$total.Range("A4").Value = 2
$total.Range("B4").Value = $oldB3
$total.Range("C4").Value = $oldC3
$total.Range("D4").Value = $oldD3

$total.Range("A3").Value = 1
$total.Range("B3").Value = $oldB2
$total.Range("C3").Value = $oldC2
$total.Range("D3").Value = $oldD2

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 2.37

# Match the index-column styling already used for A2:A3 on A4.
$a4 = $total.Range("A4")
$a4.HorizontalAlignment = -4108
$a4.VerticalAlignment = -4160
$a4.Borders.LineStyle = 1
$a4.Font.Bold = $true

# ---------------------------------------------------------------------------
# 3. Keep the first sheet ("2021-Q2") as the active / selected tab, matching
#    the unchanged bookViews/activeTab in the source workbook (Worksheets.Add
#    otherwise leaves the brand-new sheet active).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()
